# Prepend two new weekly price rows (date 44776) for "Repollo" at
# Femacal de La Calera, ahead of the existing historical rows, by
# inserting two blank rows at row 582 (pushing all following rows down
# by two) and then filling those two new rows with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 582; each Insert() shifts row 582 and
# everything below it down by one row, so doing it twice opens up a
# two-row gap at rows 582-583.
$ws.Rows.Item(582).Insert()
$ws.Rows.Item(582).Insert()

# New row 582: Primera quality
$ws.Cells.Item(582, 1).Value = 3
$ws.Cells.Item(582, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(582, 3).Value = "Coquimbo"
$ws.Cells.Item(582, 4).Value = 44776
$ws.Cells.Item(582, 5).Value = 5
$ws.Cells.Item(582, 6).Value = 100112006
$ws.Cells.Item(582, 7).Value = "Repollo"
$ws.Cells.Item(582, 8).Value = "Crespo record"
$ws.Cells.Item(582, 9).Value = "Primera"
$ws.Cells.Item(582, 10).Value = 2250
$ws.Cells.Item(582, 11).Value = 1400
$ws.Cells.Item(582, 12).Value = 1500
$ws.Cells.Item(582, 13).Value = 1449
$ws.Cells.Item(582, 14).Value = "`$/unidad"
$ws.Cells.Item(582, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(582, 16).Value = 1449
$ws.Cells.Item(582, 17).Value = 1
$ws.Cells.Item(582, 18).Value = "Hortaliza"

# New row 583: Segunda quality
$ws.Cells.Item(583, 1).Value = 3
$ws.Cells.Item(583, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(583, 3).Value = "Coquimbo"
$ws.Cells.Item(583, 4).Value = 44776
$ws.Cells.Item(583, 5).Value = 5
$ws.Cells.Item(583, 6).Value = 100112006
$ws.Cells.Item(583, 7).Value = "Repollo"
$ws.Cells.Item(583, 8).Value = "Crespo record"
$ws.Cells.Item(583, 9).Value = "Segunda"
$ws.Cells.Item(583, 10).Value = 1250
$ws.Cells.Item(583, 11).Value = 1000
$ws.Cells.Item(583, 12).Value = 1000
$ws.Cells.Item(583, 13).Value = 1000
$ws.Cells.Item(583, 14).Value = "`$/unidad"
$ws.Cells.Item(583, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(583, 16).Value = 1000
$ws.Cells.Item(583, 17).Value = 1
$ws.Cells.Item(583, 18).Value = "Hortaliza"
